$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H38").Value = 3301.7058
$ws.Range("I38").Value = 260.85715
$ws.Range("J38").Value = 5430.3
$ws.Range("K38").Value = 782.5714499999999
$ws.Range("L38").Value = 16290.9
$ws.Range("M38").Value = -410.5714499999999
$ws.Range("N38").Value = -17034.9
$ws.Range("H74").Value = 9206.25
$ws.Range("I74").Value = 5067.857
$ws.Range("K74").Value = 5067.857
$ws.Range("M74").Value = -4131.857
$ws.Range("H77").Value = 9206.25
$ws.Range("I77").Value = 5067.857
$ws.Range("K77").Value = 25339.285
$ws.Range("M77").Value = -20659.285
$ws.Range("H97").Value = 3820.75
$ws.Range("J97").Value = 3820.75
$ws.Range("L97").Value = 11462.25
$ws.Range("N97").Value = -12454.25
$ws.Range("H112").Value = 1413.409
$ws.Range("I112").Value = 1153.1428
$ws.Range("J112").Value = 1534.8667
$ws.Range("K112").Value = 3459.4284
$ws.Range("L112").Value = 4604.6001
$ws.Range("M112").Value = -2351.4284
$ws.Range("N112").Value = -6820.6001
$ws.Range("H132").Value = 1910.3243
$ws.Range("I132").Value = 1488.4333
$ws.Range("K132").Value = 4465.2999
$ws.Range("M132").Value = -1935.2999
$ws.Range("H138").Value = 6737.242
$ws.Range("J138").Value = 9567.450000000001
$ws.Range("L138").Value = 28702.35
$ws.Range("N138").Value = -38982.35000000001
$ws.Range("H141").Value = 4166.6665
$ws.Range("I141").Value = 4500
$ws.Range("J141").Value = 4000
$ws.Range("K141").Value = 13500
$ws.Range("L141").Value = 12000
$ws.Range("M141").Value = -8320
$ws.Range("N141").Value = -22360

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 2286.8794
$ws.Range("I32").Value = 2302.5818
$ws.Range("K32").Value = 2302.5818
$ws.Range("M32").Value = -2015.5818
$ws.Range("H61").Value = 2389.9302
$ws.Range("I61").Value = 2389.9302
$ws.Range("J61").Value = 0
$ws.Range("K61").Value = 2389.9302
$ws.Range("L61").Value = 0
$ws.Range("M61").Value = -2177.9302
$ws.Range("N61").Value = $null
$ws.Range("H74").Value = 3145.054
$ws.Range("I74").Value = 2479.6453
$ws.Range("J74").Value = 6583
$ws.Range("K74").Value = 2479.6453
$ws.Range("L74").Value = 6583
$ws.Range("M74").Value = -1605.6453
$ws.Range("N74").Value = -8331
$ws.Range("H77").Value = 3145.054
$ws.Range("I77").Value = 2479.6453
$ws.Range("J77").Value = 6583
$ws.Range("K77").Value = 12398.2265
$ws.Range("L77").Value = 32915
$ws.Range("M77").Value = -8030.226500000001
$ws.Range("N77").Value = -41651
$ws.Range("H102").Value = 3766.8235
$ws.Range("I102").Value = 1926
$ws.Range("K102").Value = 1926
$ws.Range("M102").Value = -304
$ws.Range("H136").Value = 2389.9302
$ws.Range("I136").Value = 2389.9302
$ws.Range("J136").Value = 0
$ws.Range("K136").Value = 7169.790599999999
$ws.Range("L136").Value = 0
$ws.Range("M136").Value = -4619.790599999999
$ws.Range("N136").Value = $null
$ws.Range("H137").Value = 150780
$ws.Range("J137").Value = 150780
$ws.Range("L137").Value = 150780
$ws.Range("N137").Value = -160980

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H22").Value = 529.8333
$ws.Range("I22").Value = 588.9
$ws.Range("J22").Value = 456
$ws.Range("K22").Value = 588.9
$ws.Range("L22").Value = 456
$ws.Range("M22").Value = -415.9
$ws.Range("N22").Value = -802
$ws.Range("H94").Value = 120359.1
$ws.Range("I94").Value = 174277.66
$ws.Range("K94").Value = 174277.66
$ws.Range("M94").Value = -173826.66
$ws.Range("H105").Value = 4404
$ws.Range("I105").Value = 4404
$ws.Range("K105").Value = 4404
$ws.Range("M105").Value = -2657

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 1787.2
$ws.Range("J58").Value = 1783.3334
$ws.Range("L58").Value = 1783.3334
$ws.Range("N58").Value = -2189.3334
$ws.Range("H105").Value = 1792.2858
$ws.Range("J105").Value = 2000
$ws.Range("L105").Value = 2000
$ws.Range("N105").Value = -5494
$ws.Range("H107").Value = 779.0625
$ws.Range("I107").Value = 805.0769
$ws.Range("K107").Value = 805.0769
$ws.Range("M107").Value = 1114.9231
$ws.Range("H133").Value = 0
$ws.Range("I133").Value = 0
$ws.Range("J133").Value = 0
$ws.Range("K133").Value = 0
$ws.Range("L133").Value = 0
$ws.Range("M133").Value = $null
$ws.Range("N133").Value = $null
$ws.Range("H136").Value = 1787.2
$ws.Range("J136").Value = 1783.3334
$ws.Range("L136").Value = 5350.0002
$ws.Range("N136").Value = -10450.0002

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H75").Value = 2737.6667
$ws.Range("I75").Value = 1606.5
$ws.Range("J75").Value = 5000
$ws.Range("K75").Value = 4819.5
$ws.Range("L75").Value = 15000
$ws.Range("M75").Value = -3821.5
$ws.Range("N75").Value = -16996
$ws.Range("H78").Value = 2737.6667
$ws.Range("I78").Value = 1606.5
$ws.Range("J78").Value = 5000
$ws.Range("K78").Value = 14458.5
$ws.Range("L78").Value = 45000
$ws.Range("M78").Value = -9466.5
$ws.Range("N78").Value = -54984

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 1482.1395
$ws.Range("I132").Value = 1252.4445
$ws.Range("J132").Value = 2663.4285
$ws.Range("K132").Value = 3757.3335
$ws.Range("L132").Value = 7990.2855
$ws.Range("M132").Value = -1227.3335
$ws.Range("N132").Value = -13050.2855

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H81").Value = 0
$ws.Range("J81").Value = 0
$ws.Range("L81").Value = 0
$ws.Range("N81").Value = $null
$ws.Range("H84").Value = 0
$ws.Range("J84").Value = 0
$ws.Range("L84").Value = 0
$ws.Range("N84").Value = $null
$ws.Range("H132").Value = 4225.8604
$ws.Range("I132").Value = 4019.7778
$ws.Range("J132").Value = 5285.7144
$ws.Range("K132").Value = 12059.3334
$ws.Range("L132").Value = 15857.1432
$ws.Range("M132").Value = -9529.3334
$ws.Range("N132").Value = -20917.1432
$ws.Range("H135").Value = 59985
$ws.Range("J135").Value = 59985
$ws.Range("L135").Value = 59985
$ws.Range("N135").Value = -70125

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H80").Value = 0
$ws.Range("I80").Value = 0
$ws.Range("K80").Value = 0
$ws.Range("M80").Value = $null
$ws.Range("H83").Value = 0
$ws.Range("I83").Value = 0
$ws.Range("K83").Value = 0
$ws.Range("M83").Value = $null
$ws.Range("H107").Value = 1702.875
$ws.Range("I107").Value = 1874.7142
$ws.Range("K107").Value = 5624.142599999999
$ws.Range("M107").Value = -3704.142599999999
$ws.Range("H122").Value = 2779.3157
$ws.Range("I122").Value = 2977.5625
$ws.Range("K122").Value = 8932.6875
$ws.Range("M122").Value = -6482.6875
$ws.Range("H132").Value = 1986.9546
$ws.Range("I132").Value = 1858.3846
$ws.Range("K132").Value = 5575.1538
$ws.Range("M132").Value = -3045.1538
$ws.Range("H136").Value = 2812.4707
$ws.Range("I136").Value = 2510.5833
$ws.Range("J136").Value = 3537
$ws.Range("K136").Value = 7531.749899999999
$ws.Range("L136").Value = 10611
$ws.Range("M136").Value = -4981.749899999999
$ws.Range("N136").Value = -15711
